$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 106
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 58
